$d = $word.ActiveDocument

# 1. Replace the title text "2.2 - Debate I" -> "Placeholder - Check Back Later".
$d.Content.Find.Execute("2.2 - Debate I", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Placeholder - Check Back Later", 2)

# 2. Remove the trailing " " and ":::" runs that follow
#    "...general edification later." — locate the exact span with Find and
#    clear it in place so the surrounding runs are left untouched.
$r = $d.Content
$found = $r.Find.Execute(" :::", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if ($found) {
    $r.Text = ""
}
